$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Split the old "Custom Turntable Belts" row (row 40) into two rows:
#    a new row 40 (390mm belt) and row 41 (210mm belt, reusing the old
#    row's content slot but with updated text/qty).
# ------------------------------------------------------------------
$ws.Rows("40").Insert()

# The insert pulled down formatting (styles) from row 39 into the new
# blank row 40 for columns that shouldn't have any value there - clear
# those back out so only A/B/D carry content, like the target file.
$ws.Range("C40").Clear()
$ws.Range("E40").Clear()
$ws.Range("F40").Clear()
$ws.Range("G40").Clear()

# New row 40: Custom Turntable Belt - ~390mm circumference variant
$ws.Range("A40").Value = "Custom Turntable Belt"
$ws.Range("B40").Value = "~390mm circumference (measured from belt)"
$ws.Range("D40").Value = 1

# Row 41 (originally row 40): Custom Turntable Belt - ~210mm circumference variant
$ws.Range("A41").Value = "Custom Turntable Belt"
$ws.Range("B41").Value = "~210mm circumference (measured from belt)"
$ws.Range("D41").Value = 1
$ws.Range("F41").Clear()
$ws.Range("G41").Clear()

# One hyperlink spans both belt rows (A40:A41), like the source workbook.
$ws.Hyperlinks.Add($ws.Range("A40:A41"), "https://www.amazon.com/gp/product/B07Y7X9FPS", "", "", "Custom Turntable Belt")
$ws.Range("A40").Style = "Hyperlink"

# ------------------------------------------------------------------
# 2) Rows 42-46 are the old rows 41-45 shifted down by the insert above;
#    their content/formulas/hyperlinks already carry over unchanged.
# ------------------------------------------------------------------

# ------------------------------------------------------------------
# 3) Append a new row 47 for "Rubber Feet".
# ------------------------------------------------------------------
$ws.Range("A47").Value = "Rubber Feet"
$ws.Range("B47").Value = "1"" diameter"
$ws.Range("C47").Value = 4
$ws.Range("D47").Value = 4
$ws.Range("E47").Value = 0
$ws.Range("G47").Value = "Amazon"
$ws.Hyperlinks.Add($ws.Range("A47"), "https://www.amazon.com/gp/product/B07VSZ2T4S")
$ws.Range("A47").Style = "Hyperlink"

# ------------------------------------------------------------------
# 4) Update the saved view state to match the edited document.
# ------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollRow = 19
$ws.Range("B51").Select()
